# Apply the cryptos-list price/volume refresh described in the commit.
# Cells that look like plain numbers (e.g. "243.65", "4.570") must stay as
# literal text (matching the original inlineStr cells), so we briefly force a
# text NumberFormat while assigning them, then restore the original (default)
# cell style by copying it from an untouched cell (B2) that still has it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style

$ws.Range('D2').Value = '26.447.03'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '1.726.63'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.65'
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4881'
$ws.Range('D7').Style = $defaultStyle
$ws.Range('E7').Value = '  +1.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2617'
$ws.Range('D8').Style = $defaultStyle
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06177'
$ws.Range('D9').Style = $defaultStyle
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '1.730.83'
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.48'
$ws.Range('D12').Style = $defaultStyle
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.570'
$ws.Range('D13').Style = $defaultStyle
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6000'
$ws.Range('D14').Style = $defaultStyle
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.24'
$ws.Range('D15').Style = $defaultStyle
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '26.470.36'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007095'
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').Value = '  +2.46%  '
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('D21').Value = '1.956.97'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.474'
$ws.Range('D22').Style = $defaultStyle
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.604'
$ws.Range('D23').Style = $defaultStyle
$ws.Range('E23').Value = '  -3.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.186'
$ws.Range('D24').Style = $defaultStyle
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.00'
$ws.Range('D25').Style = $defaultStyle
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.28'
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.401'
$ws.Range('D27').Style = $defaultStyle
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.71'
$ws.Range('D28').Style = $defaultStyle
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.715'
$ws.Range('D29').Style = $defaultStyle
$ws.Range('E29').Value = '  -4.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.956'
$ws.Range('D30').Style = $defaultStyle
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04525'
$ws.Range('D33').Style = $defaultStyle
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.616'
$ws.Range('D34').Style = $defaultStyle
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.001'
$ws.Range('D35').Style = $defaultStyle
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6233'
$ws.Range('D36').Style = $defaultStyle
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9097'
$ws.Range('D37').Style = $defaultStyle
$ws.Range('E37').Value = '  -1.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.988'
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').Value = '  -5.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.405'
$ws.Range('D39').Style = $defaultStyle
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01485'
$ws.Range('D41').Style = $defaultStyle
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.12'
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').Value = '  -4.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.426'
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3864'
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.680'
$ws.Range('D45').Style = $defaultStyle
$ws.Range('E45').Value = '  -3.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1156'
$ws.Range('D46').Style = $defaultStyle
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05362'
$ws.Range('D47').Style = $defaultStyle
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.29'
$ws.Range('D48').Style = $defaultStyle
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.706'
$ws.Range('D49').Style = $defaultStyle
$ws.Range('E49').Value = '  -1.40%  '
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.01'
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').Value = '  -0.44%  '
